$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C2").Value = 3
$ws.Range("D2").Value = 7.6
$ws.Range("E2").Value = 0.13
$ws.Range("F2").Value = 459
$ws.Range("G2").Value = 0.0001
$ws.Range("I2").Value = 1000
$ws.Range("N2").Value = "[0.002]"

$wb.Save()
